$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.035.99'
$ws.Range("E2").Value = '  +2.58%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.653.18'
$ws.Range("E3").Value = '  +3.54%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.13'
$ws.Range("E5").Value = '  +1.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.510'
$ws.Range("E6").Value = '  +1.76%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +1.61%  '
$ws.Range("E9").Value = '  +1.59%  '
$ws.Range("E10").Value = '  +4.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0867'
$ws.Range("E11").Value = '  +1.41%  '
$ws.Range("E12").Value = '  +3.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.668.55'
$ws.Range("E13").Value = '  +4.84%  '
$ws.Range("E14").Value = '  +2.17%  '
$ws.Range("E15").Value = '  +3.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.28'
$ws.Range("E16").Value = '  +2.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '239.64'
$ws.Range("E17").Value = '  +4.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '27.035.27'
$ws.Range("E18").Value = '  +2.70%  '
$ws.Range("E19").Value = '  +2.32%  '
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("E22").Value = '  +4.29%  '
$ws.Range("E23").Value = '  +2.92%  '
$ws.Range("E24").Value = '  +3.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.09'
$ws.Range("E25").Value = '  -0.27%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  +2.05%  '
$ws.Range("E28").Value = '  +1.61%  '
$ws.Range("E29").Value = '  +2.89%  '
$ws.Range("E30").Value = '  +0.69%  '
$ws.Range("E31").Value = '  +1.79%  '
$ws.Range("E32").Value = '  +3.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.517.10'
$ws.Range("E33").Value = '  +0.94%  '
$ws.Range("E34").Value = '  +5.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.60'
$ws.Range("E35").Value = '  +8.62%  '
$ws.Range("E36").Value = '  -0.22%  '
$ws.Range("E37").Value = '  +1.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.891'
$ws.Range("E38").Value = '  +8.94%  '
$ws.Range("E39").Value = '  +3.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.97'
$ws.Range("E40").Value = '  +3.21%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("E42").Value = '  +4.08%  '
$ws.Range("E43").Value = '  +8.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.793.58'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.776'
$ws.Range("E45").Value = '  +2.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.917'
$ws.Range("E46").Value = '  -2.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '89.77'
$ws.Range("E47").Value = '  +1.48%  '
$ws.Range("E48").Value = '  +0.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.52'
$ws.Range("E49").Value = '  +2.87%  '
$ws.Range("E50").Value = '  +1.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0976'
$ws.Range("E51").Value = '  +1.86%  '
